$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181, shifting existing rows 181..221 down to 182..222
$ws.Rows("181:181").Insert()

# Populate the newly inserted row 181 with the new record's data
$ws.Range("A181").Value = 10
$ws.Range("B181").Value = "Vega Modelo de Temuco"
$ws.Range("C181").Value = "La Araucanía"
$ws.Range("D181").Value = 44508
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = 100114013
$ws.Range("G181").Value = "Zanahoria"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 90
$ws.Range("K181").Value = 7000
$ws.Range("L181").Value = 7000
$ws.Range("M181").Value = 7000
$ws.Range("N181").Value = "$/saco 20 kilos"
$ws.Range("O181").Value = "Región del Maule"
$ws.Range("P181").Value = 350
$ws.Range("Q181").Value = 20
$ws.Range("R181").Value = "Hortaliza"
